# Update the results worksheet with the full set of uploaded competition
# results (time-of-day labels replaced by competitor/category ids 71-81,
# with some categories now having 2 placements instead of 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("00:50:05", "01:50:00", "01:50:00", 71, 1, 5),
    @("00:50:05", "01:50:00", "01:50:00", 71, 3, 5),
    @("01:00:00", "01:40:00", "01:40:00", 72, 2, 5),
    @("01:00:00", "01:40:00", "01:40:00", 72, 4, 5),
    @("02:00:00", "02:40:00", "02:40:00", 73, 2, 5),
    @("03:00:00", "03:40:00", "03:40:00", 74, 1, 5),
    @("04:00:00", "04:40:00", "04:40:00", 75, 2, 5),
    @("04:00:00", "04:40:00", "04:40:00", 75, 4, 5),
    @("05:00:00", "05:40:00", "05:40:00", 76, 3, 5),
    @("06:00:00", "06:40:00", "06:40:00", 77, 4, 5),
    @("07:00:00", "07:40:00", "07:40:00", 78, 1, 5),
    @("07:00:00", "07:40:00", "07:40:00", 78, 3, 5),
    @("08:00:00", "08:40:00", "08:40:00", 79, 1, 5),
    @("08:00:00", "08:40:00", "08:40:00", 79, 3, 5),
    @("09:00:00", "09:40:00", "09:40:00", 80, 1, 5),
    @("09:10:00", "09:50:00", "09:50:00", 81, 3, 5)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$ws.Range("E18").Select()
